$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" (strikeout count) values for rows 2-39, column G,
# replacing the old Strike# derived values with the regenerated K values.
$kValues = @(4,4,9,6,3,3,5,5,4,2,5,3,5,5,5,5,4,6,1,4,0,0,2,0,0,0,1,1,2,1,1,4,1,2,0,1,0,1)

for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
